# Shift the Notified Production Wind data forward by 10 days and update
# the corresponding notified production values (adding Imperial and Astro
# to the forecast portfolio).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Notified Production (MW)" values for rows 2..93 (column B),
# corresponding to timestamps shifted from day 45776/45777 to 45786/45787.
$newValues = @(1012,1011,1011,1010,959,957,955,953,862,859,855,852,794,792,790,788,846,848,851,854,979,985,991,997,1229,1240,1252,1264,1401,1407,1414,1421,1514,1518,1523,1527,1527,1527,1527,1526,1473,1472,1471,1470,1428,1426,1424,1422,1434,1432,1430,1429,1422,1418,1414,1410,1367,1364,1360,1356,1223,1219,1215,1211,1082,1075,1069,1063,923,918,912,907,830,826,822,818,796,795,794,794,769,769,768,767,714,712,710,708,625,623,620,618)

$dayOffset = 10

for ($row = 2; $row -le 97; $row++) {
    # Shift the timestamp in column A by 10 days.
    $oldDate = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 1).Value2 = $oldDate + $dayOffset

    # Update column B with the new notified production value, where
    # applicable (rows 2..93); rows 94..97 remain 0 and are left as-is.
    $idx = $row - 2
    if ($idx -lt $newValues.Length) {
        $ws.Cells.Item($row, 2).Value2 = $newValues[$idx]
    }
}
